# chore: update Sheets via scheduled runner
# Refreshes the cached market-price / profit columns (H:N) across several
# Leve-profit worksheets. Column layout per sheet:
#   H currentAveragePrice       I currentAveragePriceNQ   J currentAveragePriceHQ
#   K LevePriceNQ               L LevePriceHQ             M LeveProfitNQ
#   N LeveProfitHQ
# Some rows gain/lose individual M/N (or L/N) cells as part of the refresh
# (Excel omits cells that evaluate to blank rather than writing 0/blank
# values), so those are explicitly cleared or (re)written to match.

$wb = $excel.ActiveWorkbook

function Set-Row {
    param(
        $ws,
        [int]$row,
        $H, $I, $J, $K, $L, $M, $N
    )
    if ($null -ne $H) { $ws.Cells.Item($row, 8).Value = $H } else { $ws.Cells.Item($row, 8).ClearContents() }
    if ($null -ne $I) { $ws.Cells.Item($row, 9).Value = $I } else { $ws.Cells.Item($row, 9).ClearContents() }
    if ($null -ne $J) { $ws.Cells.Item($row, 10).Value = $J } else { $ws.Cells.Item($row, 10).ClearContents() }
    if ($null -ne $K) { $ws.Cells.Item($row, 11).Value = $K } else { $ws.Cells.Item($row, 11).ClearContents() }
    if ($null -ne $L) { $ws.Cells.Item($row, 12).Value = $L } else { $ws.Cells.Item($row, 12).ClearContents() }
    if ($null -ne $M) { $ws.Cells.Item($row, 13).Value = $M } else { $ws.Cells.Item($row, 13).ClearContents() }
    if ($null -ne $N) { $ws.Cells.Item($row, 14).Value = $N } else { $ws.Cells.Item($row, 14).ClearContents() }
}

# ---------------------------------------------------------------- ALC ----
$ws = $wb.Worksheets.Item("ALC")

Set-Row $ws 8  43.545456   27.9        200   83.69999999999999 600  55.30000000000001 -878
Set-Row $ws 52 809         809         0     2427               0   -2267              $null
Set-Row $ws 64 3032.4546   2955.6365   3070.8635 2955.6365 3070.8635 -2707.6365 -3566.8635
Set-Row $ws 67 3032.4546   2955.6365   3070.8635 2955.6365 3070.8635 -2097.6365 -4786.863499999999
Set-Row $ws 74 3406.1333   3537.75     3358.2727 3537.75    3358.2727 -2601.75  -5230.2727
Set-Row $ws 76 174277.92   373020.34   3927.2856 373020.34  3927.2856 -372705.34 -4557.2856
Set-Row $ws 77 3406.1333   3537.75     3358.2727 17688.75   16791.3635 -13008.75 -26151.3635
Set-Row $ws 79 174277.92   373020.34   3927.2856 373020.34  3927.2856 -371928.34 -6111.2856

# ---------------------------------------------------------------- ARM ----
$ws = $wb.Worksheets.Item("ARM")

Set-Row $ws 63  3264.182   2400.8572  4775  2400.8572  4775   -1714.8572 -6147
Set-Row $ws 66  3264.182   2400.8572  4775  12004.286  23875  -8572.286  -30739
Set-Row $ws 88  111201704  0          111201704 0      111201704 $null   -111202516
Set-Row $ws 91  111201704  0          111201704 0      111201704 $null   -111204512
Set-Row $ws 132 36312.734  51989.65   4958.9     155968.95 14876.7 -153438.95 -19936.7

# ---------------------------------------------------------------- BSM ----
$ws = $wb.Worksheets.Item("BSM")

Set-Row $ws 10  10000      0          10000  0        10000    $null      -10280
Set-Row $ws 99  1692.7878  1429.4762  2153.5833 1429.4762 2153.5833 68.52379999999994 -5149.5833
Set-Row $ws 105 2464.4546  2052.25    2700   2052.25    2700     -305.25    -6194

# ---------------------------------------------------------------- CRP ----
$ws = $wb.Worksheets.Item("CRP")

Set-Row $ws 62  2968.8572  2675  3038  2675   3038    -2051   -4286
Set-Row $ws 65  2968.8572  2675  3038  13375  15190   -10255  -21430
Set-Row $ws 105 490.16666  417.5 635.5 417.5  635.5   1329.5  -4129.5
Set-Row $ws 118 0          0     0     0      0       $null   $null
Set-Row $ws 137 32400      0     32400 0      32400   $null   -42600

# ---------------------------------------------------------------- CUL ----
$ws = $wb.Worksheets.Item("CUL")

Set-Row $ws 6 62500412 76923230 1523.3334 230769690 4570.0002 -230769577 -4796.0002

# ---------------------------------------------------------------- GSM ----
$ws = $wb.Worksheets.Item("GSM")

Set-Row $ws 70 2457630.2 4205937.5 9999.799999999999 4205937.5 9999.799999999999 -4205667.5 -10539.8
Set-Row $ws 73 2457630.2 4205937.5 9999.799999999999 4205937.5 9999.799999999999 -4205001.5 -11871.8
Set-Row $ws 80 4287.9165 6411      2771.4285 6411      2771.4285 -5413      -4767.4285
Set-Row $ws 83 4287.9165 6411      2771.4285 32055     13857.1425 -27063   -23841.1425

# ---------------------------------------------------------------- LTW ----
$ws = $wb.Worksheets.Item("LTW")

Set-Row $ws 12 2300 0 2300 0 2300 $null -2640
